$wb = $excel.ActiveWorkbook

# --- Sheet 1: ROW35-FE-LIFTER --- new row 27
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(27,1).Value = 45737.78956204861
$ws1.Cells.Item(27,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(27,2).Value = "0x01,0x90"
$ws1.Cells.Item(27,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws1.Cells.Item(27,4).Value = "0x01,0x82"
$ws1.Cells.Item(27,5).Value = "0xd"
$ws1.Cells.Item(27,6).Value = 400
$ws1.Cells.Item(27,7).Value = 568631262647114000000000.0
$ws1.Cells.Item(27,8).Value = 386
$ws1.Cells.Item(27,9).Value = 13

# --- Sheet 2: ROW35-MID-LIFTER --- new row 27
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(27,1).Value = 45737.63947496528
$ws2.Cells.Item(27,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(27,2).Value = "0x01,0x90"
$ws2.Cells.Item(27,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws2.Cells.Item(27,4).Value = "0x01,0x82"
$ws2.Cells.Item(27,5).Value = "0xe"
$ws2.Cells.Item(27,6).Value = 400
$ws2.Cells.Item(27,7).Value = 568631262647114000000000.0
$ws2.Cells.Item(27,8).Value = 386
$ws2.Cells.Item(27,9).Value = 14

# --- Sheet 3: ROW02-FE-LIFTER --- new row 27
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(27,1).Value = 45737.78812313658
$ws3.Cells.Item(27,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(27,2).Value = "0x01,0x90"
$ws3.Cells.Item(27,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws3.Cells.Item(27,4).Value = "0x01,0x82"
$ws3.Cells.Item(27,5).Value = "0x3"
$ws3.Cells.Item(27,6).Value = 400
$ws3.Cells.Item(27,7).Value = 568631262647114000000000.0
$ws3.Cells.Item(27,8).Value = 386
$ws3.Cells.Item(27,9).Value = 3

# --- Sheet 4: ROW02-MID-LIFTER --- new row 27
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(27,1).Value = 45737.8487365625
$ws4.Cells.Item(27,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws4.Cells.Item(27,2).Value = "0x01,0x90"
$ws4.Cells.Item(27,3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws4.Cells.Item(27,4).Value = "0x01,0x82"
$ws4.Cells.Item(27,5).Value = "0x3"
$ws4.Cells.Item(27,6).Value = 400
$ws4.Cells.Item(27,7).Value = 985046333984776000000000.0
$ws4.Cells.Item(27,8).Value = 386
$ws4.Cells.Item(27,9).Value = 3
